# Auto-generated Excel COM-interop script implementing the diff:
# - Inserts a new "Status" column (D), shifting Jan_2026..QoQ from D:H to E:I
# - Renames Nov_2025 header to Oct_2025 (now in column G) and refreshes all values
# - Updates all holding rows with refreshed Status + monthly values
# - Reorders / adds "Complete Exit" rows (incl. two new holdings) at rows 14-21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D so D:H shift right to E:I, making room for "Status"
$ws.Columns.Item(4).Insert()

# --- Header row ---
$ws.Range("A1").Value = "ISIN"
$ws.Range("B1").Value = "Stock Name"
$ws.Range("C1").Value = "Mutual Fund"
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Jan_2026"
$ws.Range("F1").Value = "Dec_2025"
$ws.Range("G1").Value = "Oct_2025"
$ws.Range("H1").Value = "MoM"
$ws.Range("I1").Value = "QoQ"

# --- Data rows (A2:I21) ---
# Row 2: HDFC Life Insurance Co Ltd
$ws.Range("A2").Value = "INE795G01014"
$ws.Range("B2").Value = "HDFC Life Insurance Co Ltd"
$ws.Range("C2").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D2").Value = "Adding Consistently"
$ws.Range("E2").Value = 9.956272999999999
$ws.Range("F2").Value = 6.706992
$ws.Range("G2").Value = 6.132316
$ws.Range("H2").Value = 3.249281
$ws.Range("I2").Value = 3.823956999999999

# Row 3: Samvardhana Motherson International Ltd
$ws.Range("A3").Value = "INE775A01035"
$ws.Range("B3").Value = "Samvardhana Motherson International Ltd"
$ws.Range("C3").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D3").Value = "Adding Consistently"
$ws.Range("E3").Value = 9.790407999999999
$ws.Range("F3").Value = 9.361610000000001
$ws.Range("G3").Value = 6.047439
$ws.Range("H3").Value = 0.4287979999999987
$ws.Range("I3").Value = 3.742969

# Row 4: HDFC Bank Limited
$ws.Range("A4").Value = "INE040A01034"
$ws.Range("B4").Value = "HDFC Bank Limited"
$ws.Range("C4").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D4").Value = "Adding Consistently"
$ws.Range("E4").Value = 9.502613
$ws.Range("F4").Value = 2.93099
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 6.571623000000001
$ws.Range("I4").Value = 9.502613

# Row 5: Kotak Mahindra Bank Limited
$ws.Range("A5").Value = "INE237A01036"
$ws.Range("B5").Value = "Kotak Mahindra Bank Limited"
$ws.Range("C5").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D5").Value = "Fresh Entry"
$ws.Range("E5").Value = 8.638477999999999
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 8.638477999999999
$ws.Range("I5").Value = 8.638477999999999

# Row 6: Coal India Ltd
$ws.Range("A6").Value = "INE522F01014"
$ws.Range("B6").Value = "Coal India Ltd"
$ws.Range("C6").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D6").Value = "Fresh Entry"
$ws.Range("E6").Value = 7.5848
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 7.5848
$ws.Range("I6").Value = 7.5848

# Row 7: ICICI Bank Limited
$ws.Range("A7").Value = "INE090A01021"
$ws.Range("B7").Value = "ICICI Bank Limited"
$ws.Range("C7").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D7").Value = "Adding Consistently"
$ws.Range("E7").Value = 6.567132
$ws.Range("F7").Value = 3.016362
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 3.55077
$ws.Range("I7").Value = 6.567132

# Row 8: Aurobindo Pharma Limited
$ws.Range("A8").Value = "INE406A01037"
$ws.Range("B8").Value = "Aurobindo Pharma Limited"
$ws.Range("C8").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D8").Value = "Adding Consistently"
$ws.Range("E8").Value = 6.310856
$ws.Range("F8").Value = 5.564083
$ws.Range("G8").Value = 5.018149
$ws.Range("H8").Value = 0.7467730000000001
$ws.Range("I8").Value = 1.292707

# Row 9: Ventive Hospitality Limited
$ws.Range("A9").Value = "INE781S01027"
$ws.Range("B9").Value = "Ventive Hospitality Limited"
$ws.Range("C9").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D9").Value = "Adding Consistently"
$ws.Range("E9").Value = 6.171473
$ws.Range("F9").Value = 5.68026
$ws.Range("G9").Value = 5.153892
$ws.Range("H9").Value = 0.4912130000000001
$ws.Range("I9").Value = 1.017581

# Row 10: Bajaj Auto Limited
$ws.Range("A10").Value = "INE917I01010"
$ws.Range("B10").Value = "Bajaj Auto Limited"
$ws.Range("C10").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D10").Value = "Adding Consistently"
$ws.Range("E10").Value = 4.906354
$ws.Range("F10").Value = 4.298993
$ws.Range("G10").Value = 3.833128
$ws.Range("H10").Value = 0.607361
$ws.Range("I10").Value = 1.073226

# Row 11: Adani Green Energy Limited
$ws.Range("A11").Value = "INE364U01010"
$ws.Range("B11").Value = "Adani Green Energy Limited"
$ws.Range("C11").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D11").Value = "Reducing Consistently"
$ws.Range("E11").Value = 4.122146
$ws.Range("F11").Value = 4.418953
$ws.Range("G11").Value = 4.649052
$ws.Range("H11").Value = -0.2968070000000003
$ws.Range("I11").Value = -0.5269060000000003

# Row 12: OSWAL PUMPS LIMITED
$ws.Range("A12").Value = "INE0BYP01024"
$ws.Range("B12").Value = "OSWAL PUMPS LIMITED"
$ws.Range("C12").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D12").Value = "Reducing Consistently"
$ws.Range("E12").Value = 1.849712
$ws.Range("F12").Value = 2.216061
$ws.Range("G12").Value = 2.859003
$ws.Range("H12").Value = -0.3663489999999998
$ws.Range("I12").Value = -1.009291

# Row 13: Capital Infra Trust InvIT
$ws.Range("A13").Value = "INE0Z8Z23013"
$ws.Range("B13").Value = "Capital Infra Trust InvIT"
$ws.Range("C13").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D13").Value = "Adding Consistently"
$ws.Range("E13").Value = 0.08158700000000001
$ws.Range("F13").Value = 0.07441300000000001
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0.007174
$ws.Range("I13").Value = 0.08158700000000001

# Row 14: Tata Power Company Limited
$ws.Range("A14").Value = "INE245A01021"
$ws.Range("B14").Value = "Tata Power Company Limited"
$ws.Range("C14").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D14").Value = "Complete Exit"
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 9.512938999999999
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = -9.512938999999999

# Row 15: Sun Pharmaceutical Industries Limited
$ws.Range("A15").Value = "INE044A01036"
$ws.Range("B15").Value = "Sun Pharmaceutical Industries Limited"
$ws.Range("C15").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D15").Value = "Complete Exit"
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 2.644174
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = -2.644174

# Row 16: Bajaj Finance Limited
$ws.Range("A16").Value = "INE296A01032"
$ws.Range("B16").Value = "Bajaj Finance Limited"
$ws.Range("C16").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D16").Value = "Complete Exit"
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1.781996
$ws.Range("G16").Value = 3.219585
$ws.Range("H16").Value = -1.781996
$ws.Range("I16").Value = -3.219585

# Row 17: Kotak Mahindra Bank Limited
$ws.Range("A17").Value = "INE237A01028"
$ws.Range("B17").Value = "Kotak Mahindra Bank Limited"
$ws.Range("C17").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D17").Value = "Complete Exit"
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 8.389315
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = -8.389315
$ws.Range("I17").Value = 0

# Row 18: Godrej Properties Limited
$ws.Range("A18").Value = "INE484J01027"
$ws.Range("B18").Value = "Godrej Properties Limited"
$ws.Range("C18").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D18").Value = "Complete Exit"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 2.567573
$ws.Range("G18").Value = 2.745639
$ws.Range("H18").Value = -2.567573
$ws.Range("I18").Value = -2.745639

# Row 19: ITC Limited
$ws.Range("A19").Value = "INE154A01025"
$ws.Range("B19").Value = "ITC Limited"
$ws.Range("C19").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D19").Value = "Complete Exit"
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 8.106636999999999
$ws.Range("G19").Value = 7.921288
$ws.Range("H19").Value = -8.106636999999999
$ws.Range("I19").Value = -7.921288

# Row 20: State Bank of India
$ws.Range("A20").Value = "INE062A01020"
$ws.Range("B20").Value = "State Bank of India"
$ws.Range("C20").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D20").Value = "Complete Exit"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 10.556393
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = -10.556393

# Row 21: DLF Limited
$ws.Range("A21").Value = "INE271C01023"
$ws.Range("B21").Value = "DLF Limited"
$ws.Range("C21").Value = "quant Dynamic Asset Allocation Fund"
$ws.Range("D21").Value = "Complete Exit"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 2.922047
$ws.Range("G21").Value = 3.011564
$ws.Range("H21").Value = -2.922047
$ws.Range("I21").Value = -3.011564

Write-Output "Edit complete"
